# Auto-generated Excel COM-interop script applying the Maduin_Profits leve-price update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2142
$ws.Range("I15").Value = 2142
$ws.Range("K15").Value = 6426
$ws.Range("M15").Value = -6257

$ws.Range("H17").Value = 746.8
$ws.Range("J17").Value = 691.0909
$ws.Range("L17").Value = 2073.2727
$ws.Range("N17").Value = -2409.2727

$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H43").Value = 5116
$ws.Range("I43").Value = 4449.25
$ws.Range("J43").Value = 6449.5
$ws.Range("K43").Value = 4449.25
$ws.Range("L43").Value = 6449.5
$ws.Range("M43").Value = -4380.25
$ws.Range("N43").Value = -6587.5

$ws.Range("H75").Value = 43000
$ws.Range("J75").Value = 43000
$ws.Range("L75").Value = 43000
$ws.Range("N75").Value = -44872

$ws.Range("H78").Value = 43000
$ws.Range("J78").Value = 43000
$ws.Range("L78").Value = 129000
$ws.Range("N78").Value = -138360

$ws.Range("H88").Value = 3249.6667
$ws.Range("I88").Value = 2999
$ws.Range("J88").Value = 3375
$ws.Range("K88").Value = 2999
$ws.Range("L88").Value = 3375
$ws.Range("M88").Value = -2593
$ws.Range("N88").Value = -4187

$ws.Range("H91").Value = 3249.6667
$ws.Range("I91").Value = 2999
$ws.Range("J91").Value = 3375
$ws.Range("K91").Value = 2999
$ws.Range("L91").Value = 3375
$ws.Range("M91").Value = -1595
$ws.Range("N91").Value = -6183

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1003
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1003
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 1003
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -1233

$ws.Range("H10").Value = 10005
$ws.Range("J10").Value = 10005
$ws.Range("L10").Value = 10005
$ws.Range("N10").Value = -10345

$ws.Range("H12").Value = 10000
$ws.Range("J12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("N12").Value = -10346

$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 2000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -2350

$ws.Range("H32").Value = 2263.3333
$ws.Range("I32").Value = 2069.8928
$ws.Range("K32").Value = 2069.8928
$ws.Range("M32").Value = -1782.8928

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H102").Value = 1665.5
$ws.Range("I102").Value = 1680.6364
$ws.Range("K102").Value = 1680.6364
$ws.Range("M102").Value = -58.63640000000009

$ws.Range("H110").Value = 1041
$ws.Range("I110").Value = 1051.25
$ws.Range("K110").Value = 1051.25
$ws.Range("M110").Value = 993.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 516.6667
$ws.Range("J8").Value = 300
$ws.Range("L8").Value = 300
$ws.Range("N8").Value = -580

$ws.Range("H11").Value = 683
$ws.Range("I11").Value = 50
$ws.Range("J11").Value = 999.5
$ws.Range("K11").Value = 50
$ws.Range("L11").Value = 999.5
$ws.Range("M11").Value = 90
$ws.Range("N11").Value = -1279.5

$ws.Range("H94").Value = 4641.357
$ws.Range("I94").Value = 4498.143
$ws.Range("K94").Value = 4498.143
$ws.Range("M94").Value = -4047.143

$ws.Range("H105").Value = 2764.2856
$ws.Range("I105").Value = 2891.6667
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2891.6667
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -1144.6667
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 19
$ws.Range("I26").Value = 19
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 19
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 268
$ws.Range("N26").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 999.8333
$ws.Range("I11").Value = 999
$ws.Range("K11").Value = 2997
$ws.Range("M11").Value = -2857

$ws.Range("H12").Value = 859.875
$ws.Range("J12").Value = 859.875
$ws.Range("L12").Value = 2579.625
$ws.Range("N12").Value = -2925.625

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H36").Value = 434
$ws.Range("I36").Value = 434
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1302
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1133
$ws.Range("N36").ClearContents()

$ws.Range("H38").Value = 187.88889
$ws.Range("I38").Value = 78.333336
$ws.Range("K38").Value = 235.000008
$ws.Range("M38").Value = 111.999992

$ws.Range("H81").Value = 1333.6666
$ws.Range("J81").Value = 1333.6666
$ws.Range("L81").Value = 4000.9998
$ws.Range("N81").Value = -6246.9998

$ws.Range("H84").Value = 1333.6666
$ws.Range("J84").Value = 1333.6666
$ws.Range("L84").Value = 12002.9994
$ws.Range("N84").Value = -23234.9994

$ws.Range("H112").Value = 46443.668
$ws.Range("I112").Value = 34000
$ws.Range("K112").Value = 102000
$ws.Range("M112").Value = -100892

$ws.Range("H115").Value = 500
$ws.Range("I115").Value = 500
$ws.Range("K115").Value = 1500
$ws.Range("M115").Value = -325

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1905.7059
$ws.Range("I102").Value = 1806.0625
$ws.Range("K102").Value = 1806.0625
$ws.Range("M102").Value = -184.0625

$ws.Range("H132").Value = 4113.3335
$ws.Range("I132").Value = 4745.4546
$ws.Range("K132").Value = 14236.3638
$ws.Range("M132").Value = -11706.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 50000
$ws.Range("J54").Value = 50000
$ws.Range("L54").Value = 50000
$ws.Range("N54").Value = -51288

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2399.5715
$ws.Range("I81").Value = 2399.5715
$ws.Range("K81").Value = 4799.143
$ws.Range("M81").Value = -3738.143

$ws.Range("H84").Value = 2399.5715
$ws.Range("I84").Value = 2399.5715
$ws.Range("K84").Value = 23995.715
$ws.Range("M84").Value = -18691.715
